# Auto update bond ETF rates
# The sheet currently has a single data/date column C (header = most
# recent date "2025/11/26"). We are inserting 4 earlier-dated columns
# (C:F) before it and shifting the original column C's data into G,
# then tweaking a handful of cells that differ from a pure "carry the
# same value forward" copy.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Duplicate column C (header + all 24 data rows, with its styling)
#    into D, E, F and G so every new column starts out identical to C.
$ws.Range("C1:C25").Copy($ws.Range("D1:D25"))
$ws.Range("C1:C25").Copy($ws.Range("E1:E25"))
$ws.Range("C1:C25").Copy($ws.Range("F1:F25"))
$ws.Range("C1:C25").Copy($ws.Range("G1:G25"))

# 2) Fix up the header dates. G1 already holds the original "2025/11/26"
#    value (copied from C1), so only C1/D1/E1/F1 need new dates.
$ws.Range("C1").Value = "'2025/11/17"
$ws.Range("D1").Value = "'2025/11/21"
$ws.Range("E1").Value = "'2025/11/24"
$ws.Range("F1").Value = "'2025/11/25"

# 3) Row 8 (科创债ETF易方达) bumps to 60 on 11/21-11/25, back to 59 on 11/26.
$ws.Range("D8").Value = 60
$ws.Range("E8").Value = 60
$ws.Range("F8").Value = 60

# 4) Rows 20, 23, 24, 25 only have a rate on the most recent date (G);
#    clear the earlier columns that were pre-filled by the column copy.
$ws.Range("C20:F20").Value = ""
$ws.Range("C23:F23").Value = ""
$ws.Range("C24:F24").Value = ""
$ws.Range("C25:F25").Value = ""

Write-Host "Bond ETF rates updated"
